# Update workbook according to the "Actualización desde MV -datos-" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 values (columns I, J, K, L) ---
$ws.Range("I74").Value = 89
$ws.Range("J74").Value = -34
$ws.Range("K74").Value = -214
$ws.Range("L74").Value = -891

# --- Append a new row 75 with the next quarterly period ---
# Force column A to be treated as text so the period label "01-04-2021"
# is stored as a shared string (like all the other period labels),
# instead of being auto-parsed into a date serial number.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()

$ws.Range("B75").Value = -376
$ws.Range("C75").Value = 2
$ws.Range("D75").Value = -378
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = -372
$ws.Range("G75").Value = -356
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = -110
$ws.Range("J75").Value = -22
$ws.Range("K75").Value = 51
$ws.Range("L75").Value = -480
$ws.Range("M75").Value = 562
$ws.Range("N75").Value = -17
